$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.58
$ws.Range("I2").Value = 7.8
$ws.Range("J2").Value = 3.95
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 3.75
$ws.Range("R2").Value = 1.37
$ws.Range("S2").Value = 3.3
$ws.Range("T2").Value = 1.93
$ws.Range("U2").Value = 1.89
$ws.Range("V2").Value = 1.15
$ws.Range("W2").Value = 2.42
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 29
$ws.Range("AB2").Value = 10.5
$ws.Range("AC2").Value = 12.5
$ws.Range("AD2").Value = 32
$ws.Range("AF2").Value = 12.5
$ws.Range("AG2").Value = 13
$ws.Range("AH2").Value = 30
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 24

# Row 3
$ws.Range("F3").Value = 2.8
$ws.Range("G3").Value = 3.25
$ws.Range("H3").Value = 2.36
$ws.Range("I3").Value = 2.7
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 3.75
$ws.Range("O3").Value = 1.28
$ws.Range("R3").Value = 1.37
$ws.Range("S3").Value = 3.1
$ws.Range("T3").Value = 1.68
$ws.Range("U3").Value = 2.16
$ws.Range("V3").Value = 1.59
$ws.Range("W3").Value = 1.44
$ws.Range("X3").Value = 19.5
$ws.Range("Y3").Value = 14
$ws.Range("Z3").Value = 20
$ws.Range("AB3").Value = 15.5
$ws.Range("AC3").Value = 10.5
$ws.Range("AD3").Value = 14.5
$ws.Range("AE3").Value = 34
$ws.Range("AF3").Value = 25
$ws.Range("AG3").Value = 16
$ws.Range("AH3").Value = 21
$ws.Range("AI3").Value = 46
$ws.Range("AJ3").Value = 60
$ws.Range("AK3").Value = 40
$ws.Range("AL3").Value = 55
$ws.Range("AM3").Value = 110
$ws.Range("AN3").Value = 34
$ws.Range("AO3").Value = 25

# Row 4
$ws.Range("G4").Value = 1.46
$ws.Range("L4").Value = 1.27
$ws.Range("M4").Value = 1.03
$ws.Range("R4").Value = 1.52
$ws.Range("V4").Value = 1.09
$ws.Range("W4").Value = 3.25

# Row 5
$ws.Range("H5").Value = 2.26
$ws.Range("I5").Value = 2.54
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.01
$ws.Range("N5").Value = 3.6
$ws.Range("O5").Value = 1.19
$ws.Range("R5").Value = 1.34
$ws.Range("S5").Value = 2.56
$ws.Range("T5").Value = 1.01
$ws.Range("U5").Value = 1.01
$ws.Range("V5").Value = 1.64
$ws.Range("W5").Value = 1.38
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 1000
$ws.Range("AD5").Value = 1000
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1000
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 1000

# Row 7
$ws.Range("H7").Value = 9.199999999999999
$ws.Range("T7").Value = 2.6

# Row 8
$ws.Range("P8").Value = 1.42

# Row 10
$ws.Range("F10").Value = 1.49

# Row 11
$ws.Range("F11").Value = 2.3
$ws.Range("G11").Value = 2.5
$ws.Range("H11").Value = 3.75
$ws.Range("Q11").Value = 2.68

# Row 12
$ws.Range("F12").Value = 1.52
$ws.Range("H12").Value = 8.4
$ws.Range("I12").Value = 9.6

# Row 13
$ws.Range("I13").Value = 11.5
$ws.Range("P13").Value = 2
$ws.Range("Q13").Value = 1.83
